$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set H31:H66 to 0 (was 1), matching the "Industries" column values
$ws.Range("H31:H66").Value = 0
